# Fixed LHS sampling to only sample across uncertainties (X) that vary
# (Ls still vary for all Ls) and rebuilt templates with
# PFLO:ALL_NO_STOPPING_DEFORESTATION_PLUR
#
# Concretely this renames worksheet "strategy_id-5008" to
# "strategy_id-5007" and adds a new worksheet "strategy_id-5009"
# (an exact duplicate of the just-renamed sheet, including the header
# row styling and the single data row) placed right after it.

$wb = $excel.ActiveWorkbook

# Rename the existing strategy_id-5008 sheet to strategy_id-5007.
$renamed = $wb.Worksheets.Item("strategy_id-5008")
$renamed.Name = "strategy_id-5007"

# Duplicate it (preserves all cell values/styles) and place the copy
# immediately after it; Excel names the copy "strategy_id-5007 (2)" by
# default, so rename it to the desired name.
$renamed.Copy([System.Reflection.Missing]::Value, $renamed)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "strategy_id-5009"

# Keep the first sheet as the active/selected tab, matching the
# original workbook's view state.
$wb.Worksheets.Item(1).Activate()
